$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '26.800.84'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'" + '  -0.89%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = "'" + '1.797.19'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "'" + '  -1.28%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').Value = "'" + '  -0.08%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = "'" + '309.55'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "'" + '  -0.45%  '
$ws.Range('E5').ClearFormats()
$ws.Range('E6').Value = "'" + '  -0.03%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').Value = "'" + '0.4392'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = "'" + '  +4.11%  '
$ws.Range('E7').ClearFormats()
$ws.Range('D8').Value = "'" + '0.3678'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = "'" + '  +0.48%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').Value = "'" + '0.07384'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "'" + '  +2.72%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').Value = "'" + '0.8547'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "'" + '  +1.85%  '
$ws.Range('E10').ClearFormats()
$ws.Range('D11').Value = "'" + '20.67'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = "'" + '  -0.19%  '
$ws.Range('E11').ClearFormats()
$ws.Range('D12').Value = "'" + '1.803.98'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "'" + '  -0.91%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').Value = "'" + '6.594'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "'" + '  -0.89%  '
$ws.Range('E13').ClearFormats()
$ws.Range('B14').Value = "'" + 'TRON'
$ws.Range('B14').ClearFormats()
$ws.Range('C14').Value = "'" + 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('C14').ClearFormats()
$ws.Range('D14').Value = "'" + '0.07063'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "'" + '  -0.18%  '
$ws.Range('E14').ClearFormats()
$ws.Range('B15').Value = "'" + 'Litecoin'
$ws.Range('B15').ClearFormats()
$ws.Range('C15').Value = "'" + 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('C15').ClearFormats()
$ws.Range('D15').Value = "'" + '91.77'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = "'" + '  +1.87%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').Value = "'" + '5.254'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = "'" + '  -0.45%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').Value = "'" + '1.002'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = "'" + '  -0.09%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').Value = "'" + '0.000008635'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'" + '  -1.22%  '
$ws.Range('E18').ClearFormats()
$ws.Range('E19').Value = "'" + '  -0.02%  '
$ws.Range('E19').ClearFormats()
$ws.Range('E20').Value = "'" + '  -1.09%  '
$ws.Range('E20').ClearFormats()
$ws.Range('D21').Value = "'" + '26.827.37'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'" + '  -1.07%  '
$ws.Range('E21').ClearFormats()
$ws.Range('D22').Value = "'" + '5.136'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'" + '  +0.21%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').Value = "'" + '10.78'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = "'" + '  -0.52%  '
$ws.Range('E23').ClearFormats()
$ws.Range('D24').Value = "'" + '1.976'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "'" + '  -0.30%  '
$ws.Range('E24').ClearFormats()
$ws.Range('D25').Value = "'" + '151.46'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "'" + '  -0.40%  '
$ws.Range('E25').ClearFormats()
$ws.Range('D26').Value = "'" + '2.192'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "'" + '  -2.99%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').Value = "'" + '18.32'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "'" + '  +0.35%  '
$ws.Range('E27').ClearFormats()
$ws.Range('D28').Value = "'" + '5.179'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "'" + '  -1.71%  '
$ws.Range('E28').ClearFormats()
$ws.Range('D29').Value = "'" + '117.30'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = "'" + '  +0.16%  '
$ws.Range('E29').ClearFormats()
$ws.Range('D30').Value = "'" + '0.08777'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "'" + '  +0.75%  '
$ws.Range('E30').ClearFormats()
$ws.Range('B31').Value = "'" + 'ImmutableX'
$ws.Range('B31').ClearFormats()
$ws.Range('C31').Value = "'" + 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C31').ClearFormats()
$ws.Range('D31').Value = "'" + '0.7373'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "'" + '  +0.22%  '
$ws.Range('E31').ClearFormats()
$ws.Range('B32').Value = "'" + 'ARBITRUM'
$ws.Range('B32').ClearFormats()
$ws.Range('C32').Value = "'" + 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C32').ClearFormats()
$ws.Range('D32').Value = "'" + '1.154'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = "'" + '  -1.82%  '
$ws.Range('E32').ClearFormats()
$ws.Range('D33').Value = "'" + '4.436'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "'" + '  +0.56%  '
$ws.Range('E33').ClearFormats()
$ws.Range('D34').Value = "'" + '2.881'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = "'" + '  -1.11%  '
$ws.Range('E34').ClearFormats()
$ws.Range('D35').Value = "'" + '0.9998'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = "'" + '  -0.05%  '
$ws.Range('E35').ClearFormats()
$ws.Range('D36').Value = "'" + '1.089'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = "'" + '  +0.05%  '
$ws.Range('E36').ClearFormats()
$ws.Range('D37').Value = "'" + '0.01958'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "'" + '  +0.56%  '
$ws.Range('E37').ClearFormats()
$ws.Range('D38').Value = "'" + '0.05168'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = "'" + '  -1.31%  '
$ws.Range('E38').ClearFormats()
$ws.Range('B39').Value = "'" + 'TheSandbox'
$ws.Range('B39').ClearFormats()
$ws.Range('C39').Value = "'" + 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C39').ClearFormats()
$ws.Range('D39').Value = "'" + '0.5206'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "'" + '  +3.52%  '
$ws.Range('E39').ClearFormats()
$ws.Range('B40').Value = "'" + 'FraxShare'
$ws.Range('B40').ClearFormats()
$ws.Range('C40').Value = "'" + 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C40').ClearFormats()
$ws.Range('D40').Value = "'" + '7.018'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'" + '  -4.23%  '
$ws.Range('E40').ClearFormats()
$ws.Range('D41').Value = "'" + '2.805'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "'" + '  -2.27%  '
$ws.Range('E41').ClearFormats()
$ws.Range('D42').Value = "'" + '0.1674'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = "'" + '  -0.64%  '
$ws.Range('E42').ClearFormats()
$ws.Range('D43').Value = "'" + '8.420'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = "'" + '  -1.44%  '
$ws.Range('E43').ClearFormats()
$ws.Range('D44').Value = "'" + '0.4939'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'" + '  +5.06%  '
$ws.Range('E44').ClearFormats()
$ws.Range('D45').Value = "'" + '1.985'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'" + '  +4.80%  '
$ws.Range('E45').ClearFormats()
$ws.Range('D46').Value = "'" + '10.38'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = "'" + '  -1.34%  '
$ws.Range('E46').ClearFormats()
$ws.Range('D47').Value = "'" + '104.24'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "'" + '  -1.73%  '
$ws.Range('E47').ClearFormats()
$ws.Range('D48').Value = "'" + '0.9995'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = "'" + '  -0.04%  '
$ws.Range('E48').ClearFormats()
$ws.Range('D49').Value = "'" + '1.660'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "'" + '  +0.93%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').Value = "'" + '0.06309'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "'" + '  -0.42%  '
$ws.Range('E50').ClearFormats()
$ws.Range('D51').Value = "'" + '0.9157'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = "'" + '  +1.90%  '
$ws.Range('E51').ClearFormats()
